$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5, 6 use the "Ano <year>" header pattern
$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)",
    "Custo Total (bilhões de R$)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $cells = $ws.Range("B1:E1")
    foreach ($cell in $cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $val -ne "") {
            $cell.Value = "Ano " + $val
        }
    }
}

# Sheet 4 uses the "Intervalo <range>" header pattern
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$cells4 = $ws4.Range("B1:E1")
foreach ($cell in $cells4) {
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $cell.Value = "Intervalo " + $val
    }
}
